$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 89.74896240234375
$ws.Range("B3").Value = 89.8494873046875
$ws.Range("B4").Value = 91.52709197998047
$ws.Range("B5").Value = 91.65299987792969
$ws.Range("B6").Value = 90.18277740478516
$ws.Range("B7").Value = 90.31535339355469
$ws.Range("B8").Value = 89.57313537597656
$ws.Range("B9").Value = 89.72521209716797
$ws.Range("B10").Value = 89.92959594726562
$ws.Range("B11").Value = 90.07334136962891
$ws.Range("B12").Value = 94.24617004394531
$ws.Range("B13").Value = 94.38392639160156
$ws.Range("B14").Value = 104.8970031738281
$ws.Range("B15").Value = 105.0306930541992
$ws.Range("B16").Value = 129.8565216064453
$ws.Range("B17").Value = 129.9865875244141
$ws.Range("B18").Value = 145.1612091064453
$ws.Range("B19").Value = 145.2893981933594
$ws.Range("B20").Value = 158.1814422607422
$ws.Range("B21").Value = 158.3042449951172
$ws.Range("B22").Value = 157.7744903564453
$ws.Range("B23").Value = 157.8909454345703
$ws.Range("B24").Value = 151.5064392089844
$ws.Range("B25").Value = 151.6171417236328
$ws.Range("B26").Value = 150.1309661865234
$ws.Range("B27").Value = 150.2380065917969
$ws.Range("B28").Value = 148.9734039306641
$ws.Range("B29").Value = 149.0795288085938
$ws.Range("B30").Value = 150.6069183349609
$ws.Range("B31").Value = 150.7148590087891
$ws.Range("B32").Value = 159.5521392822266
$ws.Range("B33").Value = 159.6636352539062
$ws.Range("B34").Value = 185.2494201660156
$ws.Range("B35").Value = 185.3649749755859
$ws.Range("B36").Value = 200.2365264892578
$ws.Range("B37").Value = 200.3559875488281
$ws.Range("B38").Value = 169.2181549072266
$ws.Range("B39").Value = 169.3409423828125
$ws.Range("B40").Value = 142.3582458496094
$ws.Range("B41").Value = 142.4840393066406
$ws.Range("B42").Value = 128.73974609375
$ws.Range("B43").Value = 128.8685760498047
$ws.Range("B44").Value = 119.0081939697266
$ws.Range("B45").Value = 119.1407012939453
$ws.Range("B46").Value = 109.4854049682617
$ws.Range("B47").Value = 109.6226806640625
$ws.Range("B48").Value = 110.6326065063477
$ws.Range("B49").Value = 110.7757186889648
